# Weekly price update: a new week's record is inserted at row 324 for the
# "Puerro" (leek) series at "Vega Modelo de Temuco", pushing the existing
# rows 324-347 down to 325-348 (the oldest record, previously row 347,
# becomes the new last row 348 - unchanged, just relocated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 324, shifting rows
# 324:347 down to 325:348.
$ws.Rows.Item(324).Insert()

# Populate the newly inserted row 324 with this week's data.
$ws.Range("A324").Value = 10
$ws.Range("B324").Value = "Vega Modelo de Temuco"
$ws.Range("C324").Value = "La Araucanía"
$ws.Range("D324").Value = 45223
$ws.Range("E324").Value = 9
$ws.Range("F324").Value = 100112005
$ws.Range("G324").Value = "Puerro"
$ws.Range("H324").Value = "Azul de Maquehue"
$ws.Range("I324").Value = "Primera"
$ws.Range("J324").Value = 40
$ws.Range("K324").Value = 8000
$ws.Range("L324").Value = 8000
$ws.Range("M324").Value = 8000
$ws.Range("N324").Value = "$/docena de paquetes"
$ws.Range("O324").Value = "Provincia de Cautín"
$ws.Range("P324").Value = 667
$ws.Range("Q324").Value = 12
$ws.Range("R324").Value = "Hortaliza"
